$d = $word.ActiveDocument

# Append a new paragraph after the existing "Hello dude" paragraph.
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">I have forgotten to say that </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>Juve</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> have been </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>loosing</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> while I working on Git</w:t></w:r>' +
  '</w:p>'

[void]$newPara.Range.InsertXML($xml)
